$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ I=0.8203074518761176; J=0.8725723693674974; M=2.535712666666667; N=7.607138; O=0.04494879354621957; P=0.05070282964779482; Q=0.1384220187606666; R=1.245798168846; S=0.03687183029880506; T=0.04424188819941292 }
  3  = @{ I=0.8203074518761176; J=0.8725723693674974; O=0.5715421877013505; P=0.6447070965264385; S=0.4688403156329966; T=0.5625535987641144 }
  4  = @{ I=0.8203074518761176; J=0.8725723693674974; M=1.538811333333333; N=4.616434; O=0.02727742533206951; P=0.03076929413956839; Q=0.08400217187533333; R=0.756019546878; S=0.022375875267891; T=0.02684843589112864 }
  5  = @{ I=0.8203074518761176; J=0.8725723693674974; M=19.206297; N=38.412594; O=0.3404565075487166; P=0.2560262755732715; Q=1.048452546933; R=6.290715281598; S=0.2792790101819299; T=0.2234014538973053 }
  6  = @{ I=0.8203074518761176; J=0.8725723693674974; K=3; L=1; M=0.8899256666666666; N=2.669777; O=0.0157750858716439; P=0.01779450411292666; Q=0.04858015221766666; R=0.437221369959; S=0.01294042049449515; T=0.01552699261553609 }
  7  = @{ G=0.011958; H=0.023916; I=0.1796925481238824; J=0.1274276306325027; M=2.535712666666667; N=7.607138; O=0.04494879354621957; P=0.05070282964779482; Q=0.030322052068; R=0.181932312408; S=0.008076963247414514; T=0.006460941448381903 }
  8  = @{ G=0.011958; H=0.023916; I=0.1796925481238824; J=0.1274276306325027; O=0.5715421877013505; P=0.6447070965264385; Q=0.385557222058; R=2.313343332348; S=0.1027018720683539; T=0.08215349776232427 }
  9  = @{ G=0.011958; H=0.023916; I=0.1796925481238824; J=0.1274276306325027; M=1.538811333333333; N=4.616434; O=0.02727742533206951; P=0.03076929413956839; Q=0.018401105924; R=0.110406635544; S=0.004901550064178509; T=0.00392085824843975 }
  10 = @{ G=0.011958; H=0.023916; I=0.1796925481238824; J=0.1274276306325027; M=19.206297; N=38.412594; O=0.3404565075487166; P=0.2560262755732715; Q=0.229668899526; R=0.9186755981039999; S=0.06117749736678668; T=0.03262482167596618 }
  11 = @{ G=0.011958; H=0.023916; I=0.1796925481238824; J=0.1274276306325027; K=3; L=1; M=0.8899256666666666; N=2.669777; O=0.0157750858716439; P=0.01779450411292666; Q=0.010641731122; R=0.063850386732; S=0.002834665377148749; T=0.002267511497390568 }
}

foreach ($rowKey in $data.Keys) {
  $cols = $data[$rowKey]
  foreach ($colKey in $cols.Keys) {
    $cellRef = "$colKey$rowKey"
    $ws.Range($cellRef).Value = $cols[$colKey]
  }
}
